{"js": "// Add hybrid bold + color (\"2C3E50\") highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific resume\n// bullet/impact paragraphs. Each target paragraph is located by its exact\n// current text, then cleared and rebuilt run-by-run so the metric\n// substrings land in their own bold+colored runs while the rest of the\n// sentence stays plain.\n\nconst HIGHLIGHT_COLOR = \"2C3E50\";\n\n// Each entry: the paragraph's ORIGINAL full text, and the ordered list of\n// [text, isMetric] segments that text should be rebuilt from.\nconst editPlans = [\n  {\n    find:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    parts: [\n      [\n        \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from \",\n        false,\n      ],\n      [\"23%\", true],\n      [\" to \", false],\n      [\"64%\", true],\n    ],\n  },\n  {\n    find:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    parts: [\n      [\"\\u2022 Achieved \", false],\n      [\"87%\", true],\n      [\" prediction accuracy for voter turnout vs. industry standard of \", false],\n      [\"71%\", true],\n      [\", reducing polling error margins from \", false],\n      [\"\\u00b14.2%\", true],\n      [\" to \", false],\n      [\"\\u00b12.1%\", true],\n    ],\n  },\n  {\n    find: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    parts: [\n      [\"\\u2022 Wrote RFP and analyzed bids from \", false],\n      [\"1,200\", true],\n      [\" vendors for research platform development\", false],\n    ],\n  },\n  {\n    find:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    parts: [\n      [\n        \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the \",\n        false,\n      ],\n      [\"$400M\", true],\n      [\" Polling Consortium Database at The Analyst Institute, now valued at \", false],\n      [\"$1B\", true],\n      [\"+\", false],\n    ],\n  },\n  {\n    find: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    parts: [\n      [\"\\u2022 Algorithm reduced mapping costs by \", false],\n      [\"73.5%\", true],\n      [\", saving campaigns and organizations \", false],\n      [\"$4.7M\", true],\n    ],\n  },\n  {\n    find: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    parts: [\n      [\"\\u2022 Achieved \", false],\n      [\"87%\", true],\n      [\" prediction accuracy for voter turnout vs. industry standard of \", false],\n      [\"71%\", true],\n    ],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map each plan to the (single) paragraph whose exact text matches it.\nfor (const plan of editPlans) {\n  const para = paragraphs.items.find((p) => p.text === plan.find);\n  if (!para) {\n    // Nothing to do if this paragraph text isn't present (already edited,\n    // or this document variant doesn't contain it).\n    continue;\n  }\n\n  // Clear the paragraph's existing content, then rebuild it piece by piece\n  // so every metric substring becomes its own bold + colored run while the\n  // rest of the text remains in plain runs.\n  para.getRange().insertText(\"\", \"Replace\");\n  await context.sync();\n\n  for (const [text, isMetric] of plan.parts) {\n    const inserted = para.insertText(text, \"End\");\n    if (isMetric) {\n      inserted.font.bold = true;\n      inserted.font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Add hybrid bold + color (\"2C3E50\") highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific resume\n# bullet/impact paragraphs. Each target paragraph is located by its exact\n# current text, then its content is cleared and rebuilt piece by piece so\n# the metric substrings land in their own bold+colored runs while the rest\n# of the sentence stays in plain runs.\n\nfunction Get-WdColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return $b * 65536 + $g * 256 + $r\n}\n\n$HighlightColor = Get-WdColor \"2C3E50\"\n$Bullet = [char]0x2022\n$CR = [char]13\n$PlusMinus = [char]0xB1\n\n$d = $word.ActiveDocument\n\n# Each plan: Find = exact current paragraph text (no trailing CR), Parts =\n# ordered list of (Text, IsMetric) pairs that should replace it.\n$editPlans = @(\n    @{\n        Find  = \"$Bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Parts = @(\n            , @(\"$Bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from \", $false)\n            , @(\"23%\", $true)\n            , @(\" to \", $false)\n            , @(\"64%\", $true)\n        )\n    },\n    @{\n        Find  = \"$Bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ${PlusMinus}4.2% to ${PlusMinus}2.1%\"\n        Parts = @(\n            , @(\"$Bullet Achieved \", $false)\n            , @(\"87%\", $true)\n            , @(\" prediction accuracy for voter turnout vs. industry standard of \", $false)\n            , @(\"71%\", $true)\n            , @(\", reducing polling error margins from \", $false)\n            , @(\"${PlusMinus}4.2%\", $true)\n            , @(\" to \", $false)\n            , @(\"${PlusMinus}2.1%\", $true)\n        )\n    },\n    @{\n        Find  = \"$Bullet Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Parts = @(\n            , @(\"$Bullet Wrote RFP and analyzed bids from \", $false)\n            , @(\"1,200\", $true)\n            , @(\" vendors for research platform development\", $false)\n        )\n    },\n    @{\n        Find  = \"$Bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Parts = @(\n            , @(\"$Bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the \", $false)\n            , @('$400M', $true)\n            , @(\" Polling Consortium Database at The Analyst Institute, now valued at \", $false)\n            , @('$1B', $true)\n            , @(\"+\", $false)\n        )\n    },\n    @{\n        Find  = \"$Bullet Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Parts = @(\n            , @(\"$Bullet Algorithm reduced mapping costs by \", $false)\n            , @(\"73.5%\", $true)\n            , @(\", saving campaigns and organizations \", $false)\n            , @('$4.7M', $true)\n        )\n    },\n    @{\n        Find  = \"$Bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Parts = @(\n            , @(\"$Bullet Achieved \", $false)\n            , @(\"87%\", $true)\n            , @(\" prediction accuracy for voter turnout vs. industry standard of \", $false)\n            , @(\"71%\", $true)\n        )\n    }\n)\n\nforeach ($plan in $editPlans) {\n    $targetText = $plan.Find + $CR\n\n    $match = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -eq $targetText) {\n            $match = $p\n            break\n        }\n    }\n    if ($match -eq $null) {\n        continue\n    }\n\n    $paraRange = $match.Range\n    $paraStart = $paraRange.Start\n    $paraEnd = $paraRange.End - 1   # exclude the paragraph mark\n\n    # Wipe the paragraph's existing content.\n    $clearRange = $d.Range($paraStart, $paraEnd)\n    $clearRange.Text = \"\"\n\n    # Rebuild it run-by-run so metric substrings get their own bold +\n    # colored run while everything else stays plain.\n    $pos = $paraStart\n    foreach ($part in $plan.Parts) {\n        $text = $part[0]\n        $isMetric = $part[1]\n\n        $insertionPoint = $d.Range($pos, $pos)\n        $insertionPoint.InsertAfter($text)\n\n        $segStart = $pos\n        $segEnd = $pos + $text.Length\n        if ($isMetric) {\n            $segRange = $d.Range($segStart, $segEnd)\n            $segRange.Font.Bold = 1\n            $segRange.Font.Color = $HighlightColor\n        }\n        $pos = $segEnd\n    }\n}\n"}
